$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "Tử Vi tọa thủ tại cung đối Tài Bạch",
    "Thiên Cơ tọa thủ tại cung đối Tài Bạch",
    "Thái Dương tọa thủ tại cung đối Tài Bạch",
    "Vũ Khúc tọa thủ tại cung đối Tài Bạch",
    "Thiên Đồng tọa thủ tại cung đối Tài Bạch",
    "Liêm Trinh tọa thủ tại cung đối Tài Bạch",
    "Thiên Phủ tọa thủ tại cung đối Tài Bạch",
    "Thái Âm tọa thủ tại cung đối Tài Bạch",
    "Tham Lang tọa thủ tại cung đối Tài Bạch",
    "Cự Môn tọa thủ tại cung đối Tài Bạch",
    "Thiên Tướng tọa thủ tại cung đối Tài Bạch",
    "Thiên Lương tọa thủ tại cung đối Tài Bạch",
    "Thất Sát tọa thủ tại cung đối Tài Bạch",
    "Phá Quân tọa thủ tại cung đối Tài Bạch",
    "Tử Vi đồng cung Thiên Cơ tại cung đối Tài Bạch",
    "Tử Vi đồng cung Thái Dương tại cung đối Tài Bạch",
    "Tử Vi đồng cung Vũ Khúc tại cung đối Tài Bạch",
    "Tử Vi đồng cung Thiên Đồng tại cung đối Tài Bạch",
    "Tử Vi đồng cung Liêm Trinh tại cung đối Tài Bạch",
    "Tử Vi đồng cung Thiên Phủ tại cung đối Tài Bạch",
    "Tử Vi đồng cung Thái Âm tại cung đối Tài Bạch",
    "Tử Vi đồng cung Tham Lang tại cung đối Tài Bạch",
    "Tử Vi đồng cung Cự Môn tại cung đối Tài Bạch",
    "Tử Vi đồng cung Thiên Tướng tại cung đối Tài Bạch",
    "Tử Vi đồng cung Thiên Lương tại cung đối Tài Bạch",
    "Tử Vi đồng cung Thất Sát tại cung đối Tài Bạch",
    "Tử Vi đồng cung Phá Quân tại cung đối Tài Bạch",
    "Thiên Cơ đồng cung Thái Dương tại cung đối Tài Bạch",
    "Thiên Cơ đồng cung Vũ Khúc tại cung đối Tài Bạch",
    "Thiên Cơ đồng cung Thiên Đồng tại cung đối Tài Bạch",
    "Thiên Cơ đồng cung Liêm Trinh tại cung đối Tài Bạch",
    "Thiên Cơ đồng cung Thiên Phủ tại cung đối Tài Bạch",
    "Thiên Cơ đồng cung Thái Âm tại cung đối Tài Bạch",
    "Thiên Cơ đồng cung Tham Lang tại cung đối Tài Bạch",
    "Thiên Cơ đồng cung Cự Môn tại cung đối Tài Bạch",
    "Thiên Cơ đồng cung Thiên Tướng tại cung đối Tài Bạch",
    "Thiên Cơ đồng cung Thiên Lương tại cung đối Tài Bạch",
    "Thiên Cơ đồng cung Thất Sát tại cung đối Tài Bạch",
    "Thiên Cơ đồng cung Phá Quân tại cung đối Tài Bạch",
    "Thái Dương đồng cung Vũ Khúc tại cung đối Tài Bạch",
    "Thái Dương đồng cung Thiên Đồng tại cung đối Tài Bạch",
    "Thái Dương đồng cung Liêm Trinh tại cung đối Tài Bạch",
    "Thái Dương đồng cung Thiên Phủ tại cung đối Tài Bạch",
    "Thái Dương đồng cung Thái Âm tại cung đối Tài Bạch",
    "Thái Dương đồng cung Tham Lang tại cung đối Tài Bạch",
    "Thái Dương đồng cung Cự Môn tại cung đối Tài Bạch",
    "Thái Dương đồng cung Thiên Tướng tại cung đối Tài Bạch",
    "Thái Dương đồng cung Thiên Lương tại cung đối Tài Bạch",
    "Thái Dương đồng cung Thất Sát tại cung đối Tài Bạch",
    "Thái Dương đồng cung Phá Quân tại cung đối Tài Bạch",
    "Vũ Khúc đồng cung Thiên Đồng tại cung đối Tài Bạch",
    "Vũ Khúc đồng cung Liêm Trinh tại cung đối Tài Bạch",
    "Vũ Khúc đồng cung Thiên Phủ tại cung đối Tài Bạch",
    "Vũ Khúc đồng cung Thái Âm tại cung đối Tài Bạch",
    "Vũ Khúc đồng cung Tham Lang tại cung đối Tài Bạch",
    "Vũ Khúc đồng cung Cự Môn tại cung đối Tài Bạch",
    "Vũ Khúc đồng cung Thiên Tướng tại cung đối Tài Bạch",
    "Vũ Khúc đồng cung Thiên Lương tại cung đối Tài Bạch",
    "Vũ Khúc đồng cung Thất Sát tại cung đối Tài Bạch",
    "Vũ Khúc đồng cung Phá Quân tại cung đối Tài Bạch",
    "Thiên Đồng đồng cung Liêm Trinh tại cung đối Tài Bạch",
    "Thiên Đồng đồng cung Thiên Phủ tại cung đối Tài Bạch",
    "Thiên Đồng đồng cung Thái Âm tại cung đối Tài Bạch",
    "Thiên Đồng đồng cung Tham Lang tại cung đối Tài Bạch",
    "Thiên Đồng đồng cung Cự Môn tại cung đối Tài Bạch",
    "Thiên Đồng đồng cung Thiên Tướng tại cung đối Tài Bạch",
    "Thiên Đồng đồng cung Thiên Lương tại cung đối Tài Bạch",
    "Thiên Đồng đồng cung Thất Sát tại cung đối Tài Bạch",
    "Thiên Đồng đồng cung Phá Quân tại cung đối Tài Bạch",
    "Liêm Trinh đồng cung Thiên Phủ tại cung đối Tài Bạch",
    "Liêm Trinh đồng cung Thái Âm tại cung đối Tài Bạch",
    "Liêm Trinh đồng cung Tham Lang tại cung đối Tài Bạch",
    "Liêm Trinh đồng cung Cự Môn tại cung đối Tài Bạch",
    "Liêm Trinh đồng cung Thiên Tướng tại cung đối Tài Bạch",
    "Liêm Trinh đồng cung Thiên Lương tại cung đối Tài Bạch",
    "Liêm Trinh đồng cung Thất Sát tại cung đối Tài Bạch",
    "Liêm Trinh đồng cung Phá Quân tại cung đối Tài Bạch",
    "Thiên Phủ đồng cung Thái Âm tại cung đối Tài Bạch",
    "Thiên Phủ đồng cung Tham Lang tại cung đối Tài Bạch",
    "Thiên Phủ đồng cung Cự Môn tại cung đối Tài Bạch",
    "Thiên Phủ đồng cung Thiên Tướng tại cung đối Tài Bạch",
    "Thiên Phủ đồng cung Thiên Lương tại cung đối Tài Bạch",
    "Thiên Phủ đồng cung Thất Sát tại cung đối Tài Bạch",
    "Thiên Phủ đồng cung Phá Quân tại cung đối Tài Bạch",
    "Thái Âm đồng cung Tham Lang tại cung đối Tài Bạch",
    "Thái Âm đồng cung Cự Môn tại cung đối Tài Bạch",
    "Thái Âm đồng cung Thiên Tướng tại cung đối Tài Bạch",
    "Thái Âm đồng cung Thiên Lương tại cung đối Tài Bạch",
    "Thái Âm đồng cung Thất Sát tại cung đối Tài Bạch",
    "Thái Âm đồng cung Phá Quân tại cung đối Tài Bạch",
    "Tham Lang đồng cung Cự Môn tại cung đối Tài Bạch",
    "Tham Lang đồng cung Thiên Tướng tại cung đối Tài Bạch",
    "Tham Lang đồng cung Thiên Lương tại cung đối Tài Bạch",
    "Tham Lang đồng cung Thất Sát tại cung đối Tài Bạch",
    "Tham Lang đồng cung Phá Quân tại cung đối Tài Bạch",
    "Cự Môn đồng cung Thiên Tướng tại cung đối Tài Bạch",
    "Cự Môn đồng cung Thiên Lương tại cung đối Tài Bạch",
    "Cự Môn đồng cung Thất Sát tại cung đối Tài Bạch",
    "Cự Môn đồng cung Phá Quân tại cung đối Tài Bạch",
    "Thiên Tướng đồng cung Thiên Lương tại cung đối Tài Bạch",
    "Thiên Tướng đồng cung Thất Sát tại cung đối Tài Bạch",
    "Thiên Tướng đồng cung Phá Quân tại cung đối Tài Bạch",
    "Thiên Lương đồng cung Thất Sát tại cung đối Tài Bạch",
    "Thiên Lương đồng cung Phá Quân tại cung đối Tài Bạch",
    "Thất Sát đồng cung Phá Quân tại cung đối Tài Bạch"
)

$startRow = 4324
for ($i = 0; $i -lt $values.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value() = $values[$i]
    $ws.Cells.Item($r, 2).Value() = $values[$i]
}

$null = $ws.Range("B4338:B4428").Select()

